$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 15.98786064435238
$ws.Range("C2").Value = 9.961696387662483
$ws.Range("D2").Value = 5.995628620523291
$ws.Range("E2").Value = 11.65046834840988
$ws.Range("G2").Value = 3.65322623281584
$ws.Range("I2").Value = 24.37985028251085
$ws.Range("L2").Value = 10.11610441137285
$ws.Range("M2").Value = 15.12902099349913
$ws.Range("O2").Value = 25.07308721050791
$ws.Range("B3").Value = 15.45479023350826
$ws.Range("C3").Value = 9.599675317404161
$ws.Range("D3").Value = 5.879040245493583
$ws.Range("E3").Value = 11.68528386554583
$ws.Range("G3").Value = 3.655718096751454
$ws.Range("I3").Value = 24.51574055177545
$ws.Range("L3").Value = 10.12597274198202
$ws.Range("M3").Value = 15.01538796541074
$ws.Range("O3").Value = 25.15265549392466
$ws.Range("B4").Value = 15.11967334874031
$ws.Range("C4").Value = 9.368705870796774
$ws.Range("D4").Value = 5.808098790139577
$ws.Range("E4").Value = 11.70780385855273
$ws.Range("G4").Value = 3.657329004553973
$ws.Range("I4").Value = 24.60519621182906
$ws.Range("L4").Value = 10.13349807176303
$ws.Range("M4").Value = 14.94728454412408
$ws.Range("O4").Value = 25.20873593735106
$ws.Range("B5").Value = 14.98135023465693
$ws.Range("C5").Value = 9.272490457102894
$ws.Range("D5").Value = 5.779395672162665
$ws.Range("E5").Value = 11.71726924679138
$ws.Range("G5").Value = 3.658005871253338
$ws.Range("I5").Value = 24.64316043027555
$ws.Range("L5").Value = 10.13693366897246
$ws.Range("M5").Value = 14.9199737980571
$ws.Range("O5").Value = 25.23339807946186
$ws.Range("B6").Value = 14.95828176506716
$ws.Range("C6").Value = 9.25639038748761
$ws.Range("D6").Value = 5.774643282633908
$ws.Range("E6").Value = 11.71885840572322
$ws.Range("G6").Value = 3.658119499070779
$ws.Range("I6").Value = 24.6495554901809
$ws.Range("L6").Value = 10.13752643902253
$ws.Range("M6").Value = 14.91546620831558
$ws.Range("O6").Value = 25.23760222331737
$ws.Range("B7").Value = 15.11781472020784
$ws.Range("C7").Value = 9.367416624508513
$ws.Range("D7").Value = 5.807710796893043
$ws.Range("E7").Value = 11.707930343594
$ws.Range("G7").Value = 3.657338050291463
$ws.Range("I7").Value = 24.60570209943256
$ws.Range("L7").Value = 10.13354291117073
$ws.Range("M7").Value = 14.94691440310883
$ws.Range("O7").Value = 25.20906122625679
$ws.Range("B8").Value = 15.80579543965894
$ws.Range("C8").Value = 9.8387258132861
$ws.Range("D8").Value = 5.955321570765614
$ws.Range("E8").Value = 11.6622360019664
$ws.Range("G8").Value = 3.654068678757163
$ws.Range("I8").Value = 24.4254542392413
$ws.Range("L8").Value = 10.11920290938416
$ws.Range("M8").Value = 15.08950794929404
$ws.Range("O8").Value = 25.09901811689219
$ws.Range("B9").Value = 17.08491650633286
$ws.Range("C9").Value = 10.69056312661674
$ws.Range("D9").Value = 6.247986451092766
$ws.Range("E9").Value = 11.58166040827935
$ws.Range("G9").Value = 3.64829622630393
$ws.Range("I9").Value = 24.11988249371866
$ws.Range("L9").Value = 10.10270073093371
$ws.Range("M9").Value = 15.38134440892687
$ws.Range("O9").Value = 24.94088975133822
$ws.Range("B10").Value = 17.97246992059128
$ws.Range("C10").Value = 11.26829737487886
$ws.Range("D10").Value = 6.462465827408871
$ws.Range("E10").Value = 11.52791187538058
$ws.Range("G10").Value = 3.644440312957157
$ws.Range("I10").Value = 23.9247705782477
$ws.Range("L10").Value = 10.09763728174643
$ws.Range("M10").Value = 15.60179783037752
$ws.Range("O10").Value = 24.86030190007125
$ws.Range("B11").Value = 18.3632216515337
$ws.Range("C11").Value = 11.52000875627908
$ws.Range("D11").Value = 6.559420552482326
$ws.Range("E11").Value = 11.50463205135592
$ws.Range("G11").Value = 3.642768863422949
$ws.Range("I11").Value = 23.84243492682277
$ws.Range("L11").Value = 10.09686075700018
$ws.Range("M11").Value = 15.7030932262767
$ws.Range("O11").Value = 24.83145132074812
$ws.Range("B12").Value = 18.50920040160874
$ws.Range("C12").Value = 11.61368277814033
$ws.Range("D12").Value = 6.596007443696447
$ws.Range("E12").Value = 11.49598404196276
$ws.Range("G12").Value = 3.642147740569034
$ws.Range("I12").Value = 23.81218391742416
$ws.Range("L12").Value = 10.09678558934846
$ws.Range("M12").Value = 15.74157027951695
$ws.Range("O12").Value = 24.82165527337303
$ws.Range("B13").Value = 18.4778515662171
$ws.Range("C13").Value = 11.59358218089214
$ws.Range("D13").Value = 6.58813409242177
$ws.Range("E13").Value = 11.49783910588814
$ws.Range("G13").Value = 3.642280985835169
$ws.Range("I13").Value = 23.81865768400059
$ws.Range("L13").Value = 10.09679205469718
$ws.Range("M13").Value = 15.73327870231082
$ws.Range("O13").Value = 24.823714722364
$ws.Range("B14").Value = 18.37527194482076
$ws.Range("C14").Value = 11.52774856094816
$ws.Range("D14").Value = 6.562433357918386
$ws.Range("E14").Value = 11.50391722083731
$ws.Range("G14").Value = 3.642717526761683
$ws.Range("I14").Value = 23.83992754043042
$ws.Range("L14").Value = 10.09685019028582
$ws.Range("M14").Value = 15.70625653022305
$ws.Range("O14").Value = 24.83062273359638
$ws.Range("B15").Value = 18.31217640404121
$ws.Range("C15").Value = 11.4872082374511
$ws.Range("D15").Value = 6.546673138221331
$ws.Range("E15").Value = 11.50766203916986
$ws.Range("G15").Value = 3.642986457993187
$ws.Range("I15").Value = 23.85307688616279
$ws.Range("L15").Value = 10.09691428389443
$ws.Range("M15").Value = 15.68971933699405
$ws.Range("O15").Value = 24.83500129503352
$ws.Range("B16").Value = 17.94666109771756
$ws.Range("C16").Value = 11.25161985239027
$ws.Range("D16").Value = 6.456113839469967
$ws.Range("E16").Value = 11.52945675601358
$ws.Range("G16").Value = 3.644551203381863
$ws.Range("I16").Value = 23.93028099176285
$ws.Range("L16").Value = 10.09771869798338
$ws.Range("M16").Value = 15.59519599373305
$ws.Range("O16").Value = 24.86234502811605
$ws.Range("B17").Value = 17.71900880308203
$ws.Range("C17").Value = 11.10421468458124
$ws.Range("D17").Value = 6.400372800471558
$ws.Range("E17").Value = 11.54312638144559
$ws.Range("G17").Value = 3.645532241553826
$ws.Range("I17").Value = 23.97929098192538
$ws.Range("L17").Value = 10.09860279172828
$ws.Range("M17").Value = 15.53744880627023
$ws.Range("O17").Value = 24.88112416453534
$ws.Range("B18").Value = 17.58685180394835
$ws.Range("C18").Value = 11.01838844032548
$ws.Range("D18").Value = 6.368256572065334
$ws.Range("E18").Value = 11.55109902230257
$ws.Range("G18").Value = 3.646104289508453
$ws.Range("I18").Value = 24.00808441404795
$ws.Range("L18").Value = 10.09925506032114
$ws.Range("M18").Value = 15.5043308167658
$ws.Range("O18").Value = 24.89266014120671
$ws.Range("B19").Value = 17.54190054494454
$ws.Range("C19").Value = 10.98915160949851
$ws.Range("D19").Value = 6.357374284820456
$ws.Range("E19").Value = 11.55381737827006
$ws.Range("G19").Value = 3.646299313380195
$ws.Range("I19").Value = 24.01793701567855
$ws.Range("L19").Value = 10.09950061763652
$ws.Range("M19").Value = 15.49313504496849
$ws.Range("O19").Value = 24.89669201792774
$ws.Range("B20").Value = 17.74336969876295
$ws.Range("C20").Value = 11.12001449281771
$ws.Range("D20").Value = 6.406312562362876
$ws.Range("E20").Value = 11.5416598233911
$ws.Range("G20").Value = 3.64542700354186
$ws.Range("I20").Value = 23.97401122352495
$ws.Range("L20").Value = 10.09849380432254
$ws.Range("M20").Value = 15.54358628062792
$ws.Range("O20").Value = 24.87904900732578
$ws.Range("B21").Value = 18.40545697961134
$ws.Range("C21").Value = 11.54713044041425
$ws.Range("D21").Value = 6.569986059213504
$ws.Range("E21").Value = 11.50212738996787
$ws.Range("G21").Value = 3.642588983846531
$ws.Range("I21").Value = 23.83365485550411
$ws.Range("L21").Value = 10.09682717976102
$ws.Range("M21").Value = 15.71419058793894
$ws.Range("O21").Value = 24.82856299579527
$ws.Range("B22").Value = 18.82652116487395
$ws.Range("C22").Value = 11.81667784662421
$ws.Range("D22").Value = 6.676194612325599
$ws.Range("E22").Value = 11.47726693760173
$ws.Range("G22").Value = 3.640803034276029
$ws.Range("I22").Value = 23.74733334368726
$ws.Range("L22").Value = 10.0970134041761
$ws.Range("M22").Value = 15.8263698313622
$ws.Range("O22").Value = 24.80215079092499
$ws.Range("B23").Value = 18.60289341000153
$ws.Range("C23").Value = 11.67370736768825
$ws.Range("D23").Value = 6.619591325554238
$ws.Range("E23").Value = 11.49044635662605
$ws.Range("G23").Value = 3.641749948995544
$ws.Range("I23").Value = 23.7929083513909
$ws.Range("L23").Value = 10.0967975526048
$ws.Range("M23").Value = 15.76644426461971
$ws.Range("O23").Value = 24.81564324872356
$ws.Range("B24").Value = 17.73236010501055
$ws.Range("C24").Value = 11.11287476319043
$ws.Range("D24").Value = 6.403627410921155
$ws.Range("E24").Value = 11.54232249977918
$ws.Range("G24").Value = 3.645474556609702
$ws.Range("I24").Value = 23.97639628082507
$ws.Range("L24").Value = 10.09854262895807
$ws.Range("M24").Value = 15.54081127214011
$ws.Range("O24").Value = 24.87998488248685
$ws.Range("B25").Value = 16.74742625236137
$ws.Range("C25").Value = 10.46831418878282
$ws.Range("D25").Value = 6.168732574883115
$ws.Range("E25").Value = 11.60249714509166
$ws.Range("G25").Value = 3.649789889190282
$ws.Range("I25").Value = 24.19740069582175
$ws.Range("L25").Value = 10.1059231596572
$ws.Range("M25").Value = 15.30123074996484
$ws.Range("O25").Value = 24.97744740153568
